$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap G1 and H1 cell values (header labels "Data ważności" / "Data produkcji")
$gVal = $ws.Range("G1").Value2
$hVal = $ws.Range("H1").Value2
$ws.Range("G1").Value2 = $hVal
$ws.Range("H1").Value2 = $gVal

# Update the active selection to H2
$ws.Range("H2").Select()
